# Case_4_235 (380 kV case): slack bus voltage setpoint changed from 1.05 to
# 1.02 pu, which changes the power-flow solution for all bus voltage
# magnitudes. Columns B:F and I:N (rows 2-25) are updated to the recalculated
# per-unit voltages; column G (slack reference, always 1) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: columns B:F (rows 2-25)
$block1 = New-Object 'double[,]' 24,5
$block1[0,0] = 1.02; $block1[0,1] = 1.044946558196912; $block1[0,2] = 1.048112601328319; $block1[0,3] = 1.048611548717642; $block1[0,4] = 1.053766071812368
$block1[1,0] = 1.02; $block1[1,1] = 1.04619719095088; $block1[1,2] = 1.049111611676476; $block1[1,3] = 1.049811845816773; $block1[1,4] = 1.055228126875075
$block1[2,0] = 1.02; $block1[2,1] = 1.04700538761458; $block1[2,2] = 1.049757099156007; $block1[2,3] = 1.050587822955968; $block1[2,4] = 1.056173549686864
$block1[3,0] = 1.02; $block1[3,1] = 1.047344906979854; $block1[3,2] = 1.050028239804374; $block1[3,3] = 1.050913880566358; $block1[3,4] = 1.056570861047811
$block1[4,0] = 1.02; $block1[4,1] = 1.047401899435292; $block1[4,2] = 1.050073752550088; $block1[4,3] = 1.050968617600307; $block1[4,4] = 1.0566375630555
$block1[5,0] = 1.02; $block1[5,1] = 1.047009925252615; $block1[5,2] = 1.049760723020939; $block1[5,3] = 1.050592180390805; $block1[5,4] = 1.056178859141672
$block1[6,0] = 1.02; $block1[6,1] = 1.045369434450585; $block1[6,2] = 1.048450417500197; $block1[6,3] = 1.049017340701103; $block1[6,4] = 1.054260312585511
$block1[7,0] = 1.02; $block1[7,1] = 1.042470487156607; $block1[7,2] = 1.046134186705304; $block1[7,3] = 1.046236798595462; $block1[7,4] = 1.050874574476566
$block1[8,0] = 1.02; $block1[8,1] = 1.040532094807151; $block1[8,2] = 1.044584960535302; $block1[8,3] = 1.044379214896562; $block1[8,4] = 1.048613720758785
$block1[9,0] = 1.02; $block1[9,1] = 1.039691326353508; $block1[9,2] = 1.043912888899826; $block1[9,3] = 1.04357389103668; $block1[9,4] = 1.047633797151841
$block1[10,0] = 1.02; $block1[10,1] = 1.039378807425609; $block1[10,2] = 1.043663061060415; $block1[10,3] = 1.043274606827554; $block1[10,4] = 1.047269659479148
$block1[11,0] = 1.02; $block1[11,1] = 1.039445853836269; $block1[11,2] = 1.04371665865327; $block1[11,3] = 1.043338811209641; $block1[11,4] = 1.047347775154593
$block1[12,0] = 1.02; $block1[12,1] = 1.039665497960092; $block1[12,2] = 1.043892241957853; $block1[12,3] = 1.043549155216522; $block1[12,4] = 1.04760370052336
$block1[13,0] = 1.02; $block1[13,1] = 1.039800798640569; $block1[13,2] = 1.044000399285061; $block1[13,3] = 1.043678734958284; $block1[13,4] = 1.047761364468401
$block1[14,0] = 1.02; $block1[14,1] = 1.04058786325001; $block1[14,2] = 1.044629537178329; $block1[14,3] = 1.044432640632119; $block1[14,4] = 1.048678734316043
$block1[15,0] = 1.02; $block1[15,1] = 1.041081181250008; $block1[15,2] = 1.045023842757588; $block1[15,3] = 1.044905281013992; $block1[15,4] = 1.049253914966243
$block1[16,0] = 1.02; $block1[16,1] = 1.041368787486297; $block1[16,2] = 1.045253714484456; $block1[16,3] = 1.045180870086031; $block1[16,4] = 1.049589316006085
$block1[17,0] = 1.02; $block1[17,1] = 1.04146683055595; $block1[17,2] = 1.045332074538289; $block1[17,3] = 1.045274823007508; $block1[17,4] = 1.04970366366168
$block1[18,0] = 1.02; $block1[18,1] = 1.041028267188801; $block1[18,2] = 1.04498154994211; $block1[18,3] = 1.044854580935143; $block1[18,4] = 1.049192213099958
$block1[19,0] = 1.02; $block1[19,1] = 1.03960082435578; $block1[19,2] = 1.043840542336093; $block1[19,3] = 1.043487218352983; $block1[19,4] = 1.047528341084202
$block1[20,0] = 1.02; $block1[20,1] = 1.038702058539046; $block1[20,2] = 1.043122041715458; $block1[20,3] = 1.042626627103284; $block1[20,4] = 1.046481325301344
$block1[21,0] = 1.02; $block1[21,1] = 1.039178633908872; $block1[21,2] = 1.04350303820821; $block1[21,3] = 1.043082927247878; $block1[21,4] = 1.047036453088422
$block1[22,0] = 1.02; $block1[22,1] = 1.041052177200934; $block1[22,2] = 1.045000660615752; $block1[22,3] = 1.044877490408878; $block1[22,4] = 1.049220093801591
$block1[23,0] = 1.02; $block1[23,1] = 1.043220931327096; $block1[23,2] = 1.046733870060642; $block1[23,3] = 1.046956305444677; $block1[23,4] = 1.051750496364691

# Block 2: columns I:N (rows 2-25)
$block2 = New-Object 'double[,]' 24,6
$block2[0,0] = 1.045997656216752; $block2[0,1] = 1.050008921586046; $block2[0,2] = 1.050873354813739; $block2[0,3] = 1.051370910096921; $block2[0,4] = 1.056511135245837; $block2[0,5] = 1.020279421341814
$block2[1,0] = 1.046465442010958; $block2[1,1] = 1.050905528185158; $block2[1,2] = 1.051683860775927; $block2[1,3] = 1.052382283382463; $block2[1,4] = 1.057784639903008; $block2[1,5] = 1.020602559376399
$block2[2,0] = 1.046766292966245; $block2[2,1] = 1.051484150606172; $block2[2,2] = 1.052206764253219; $block2[2,3] = 1.053035447810978; $block2[2,4] = 1.058607545703716; $block2[2,5] = 1.02081052835659
$block2[3,0] = 1.046892332445008; $block2[3,1] = 1.051727037134713; $block2[3,2] = 1.052426224714892; $block2[3,3] = 1.053309739088548; $block2[3,4] = 1.058953227685894; $block2[3,5] = 1.020897690594874
$block2[4,0] = 1.046913469377255; $block2[4,1] = 1.051767797451206; $block2[4,2] = 1.052463051574066; $block2[4,3] = 1.053355776351899; $block2[4,4] = 1.059011253650648; $block2[4,5] = 1.020912309833656
$block2[5,0] = 1.046767978831902; $block2[5,1] = 1.051487397503803; $block2[5,2] = 1.052209698138223; $block2[5,3] = 1.053039114075767; $block2[5,4] = 1.05861216576917; $block2[5,5] = 1.020811694074288
$block2[6,0] = 1.046156128152973; $block2[6,1] = 1.05031225496768; $block2[6,2] = 1.051147591117724; $block2[6,3] = 1.051712971495505; $block2[6,4] = 1.056941760248066; $block2[6,5] = 1.020388860414872
$block2[7,0] = 1.045063817995636; $block2[7,1] = 1.048229566911073; $block2[7,2] = 1.049264056036915; $block2[7,3] = 1.049366340290079; $block2[7,4] = 1.053989380935107; $block2[7,5] = 1.0196351320679
$block2[8,0] = 1.044325988896067; $block2[8,1] = 1.046832903308447; $block2[8,2] = 1.048000163887194; $block2[8,3] = 1.047795140254658; $block2[8,4] = 1.052014849184872; $block2[8,5] = 1.019126778122367
$block2[9,0] = 1.044004194019674; $block2[9,1] = 1.046226145963908; $block2[9,2] = 1.047450903820174; $block2[9,3] = 1.047113142154535; $block2[9,4] = 1.051158301395514; $block2[9,5] = 1.018905250049665
$block2[10,0] = 1.04388431584649; $block2[10,1] = 1.046000466299448; $block2[10,2] = 1.04724658260117; $block2[10,3] = 1.046859564769043; $block2[10,4] = 1.05083990102202; $block2[10,5] = 1.018822752010188
$block2[11,0] = 1.04391004595995; $block2[11,1] = 1.046048889121282; $block2[11,2] = 1.047290423893665; $block2[11,3] = 1.046913969485856; $block2[11,4] = 1.050908209951805; $block2[11,5] = 1.018840457761805
$block2[12,0] = 1.043994291989261; $block2[12,1] = 1.046207497413718; $block2[12,2] = 1.047434020734706; $block2[12,3] = 1.047092186550471; $block2[12,4] = 1.051131987257844; $block2[12,5] = 1.01889843508091
$block2[13,0] = 1.044046152406974; $block2[13,1] = 1.046305180963502; $block2[13,2] = 1.047522455432101; $block2[13,3] = 1.047201958336448; $block2[13,4] = 1.051269831818969; $block2[13,5] = 1.018934128606701
$block2[14,0] = 1.04434729650758; $block2[14,1] = 1.046873129504492; $block2[14,2] = 1.048036574349923; $block2[14,3] = 1.047840366962981; $block2[14,4] = 1.052071662039879; $block2[14,5] = 1.019141450455566
$block2[15,0] = 1.044535576258802; $block2[15,1] = 1.047228852573967; $block2[15,2] = 1.048358533533036; $block2[15,3] = 1.048240377218271; $block2[15,4] = 1.052574206653802; $block2[15,5] = 1.019271120315311
$block2[16,0] = 1.044645173918123; $block2[16,1] = 1.0474361477903; $block2[16,2] = 1.048546135429713; $block2[16,3] = 1.048473536707789; $block2[16,4] = 1.05286718204807; $block2[16,5] = 1.019346618832617
$block2[17,0] = 1.044682506173553; $block2[17,1] = 1.047506797637759; $block2[17,2] = 1.048610070449849; $block2[17,3] = 1.04855301109089; $block2[17,4] = 1.052967053749612; $block2[17,5] = 1.019372338891729
$block2[18,0] = 1.044515398678223; $block2[18,1] = 1.047190706744864; $block2[18,2] = 1.048324010185706; $block2[18,3] = 1.04819747645338; $block2[18,4] = 1.052520303969224; $block2[18,5] = 1.019257222013573
$block2[19,0] = 1.043969493291347; $block2[19,1] = 1.046160799617109; $block2[19,2] = 1.047391743394287; $block2[19,3] = 1.04703971307958; $block2[19,4] = 1.051066097099988; $block2[19,5] = 1.018881368085281
$block2[20,0] = 1.043624239295009; $block2[20,1] = 1.04551150144746; $block2[20,2] = 1.04680384387835; $block2[20,3] = 1.046310316044852; $block2[20,4] = 1.050150385907799; $block2[20,5] = 1.018643822881007
$block2[21,0] = 1.043807457325605; $block2[21,1] = 1.045855874214717; $block2[21,2] = 1.047115667123532; $block2[21,3] = 1.046697123380058; $block2[21,4] = 1.05063595572426; $block2[21,5] = 1.018769867165756
$block2[22,0] = 1.044524516746452; $block2[22,1] = 1.047207943796234; $block2[22,2] = 1.048339610391718; $block2[22,3] = 1.048216861955732; $block2[22,4] = 1.052544660736459; $block2[22,5] = 1.019263502477195
$block2[23,0] = 1.045347894405272; $block2[23,1] = 1.048769424661719; $block2[23,2] = 1.049752429164505; $block2[23,3] = 1.049974181507218; $block2[23,4] = 1.054753727300543; $block2[23,5] = 1.019831018956487

$ws.Range("B2:F25").Value = $block1
$ws.Range("I2:N25").Value = $block2

Write-Output "Updated vm_pu.xlsx bus voltages for 380 kV case (slack = 1.02 pu)"
